# Auto-generated data refresh for the Leve profit-tracker workbook.
# Each worksheet (one per crafting job) stores scraped market-board
# prices in columns H:N; this script overwrites the scraped cells with
# the latest pull, clearing any cell that the refresh left blank.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = ""
$ws.Range("M2").Value = 13
$ws.Range("H18").Value = 366.66666
$ws.Range("I18").Value = 366.66666
$ws.Range("K18").Value = 366.66666
$ws.Range("M18").Value = -82.66665999999998
$ws.Range("H40").Value = 1206.6
$ws.Range("J40").Value = 1175
$ws.Range("L40").Value = 1175
$ws.Range("N40").Value = -1525
$ws.Range("H58").Value = 181.5
$ws.Range("I58").Value = 193.125
$ws.Range("J58").Value = 135
$ws.Range("K58").Value = 579.375
$ws.Range("L58").Value = 405
$ws.Range("M58").Value = -429.375
$ws.Range("N58").Value = -705
$ws.Range("H116").Value = 14930.8
$ws.Range("J116").Value = 17453.223
$ws.Range("L116").Value = 17453.223
$ws.Range("N116").Value = -24337.223
$ws.Range("H138").Value = 4438.4443
$ws.Range("J138").Value = 4724.1587
$ws.Range("L138").Value = 14172.4761
$ws.Range("N138").Value = -24452.4761
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2305.8
$ws.Range("I2").Value = 2284.2222
$ws.Range("K2").Value = 2284.2222
$ws.Range("M2").Value = -2171.2222
$ws.Range("H32").Value = 2539.25
$ws.Range("I32").Value = 1906.0892
$ws.Range("J32").Value = 6971.375
$ws.Range("K32").Value = 1906.0892
$ws.Range("L32").Value = 6971.375
$ws.Range("M32").Value = -1619.0892
$ws.Range("N32").Value = -7545.375
$ws.Range("H45").Value = 21873.75
$ws.Range("J45").Value = 18201.666
$ws.Range("L45").Value = 18201.666
$ws.Range("N45").Value = -18955.666
$ws.Range("H61").Value = 3983.8667
$ws.Range("I61").Value = 3519.923
$ws.Range("K61").Value = 3519.923
$ws.Range("M61").Value = -3307.923
$ws.Range("H116").Value = 2305.8
$ws.Range("I116").Value = 2284.2222
$ws.Range("K116").Value = 2284.2222
$ws.Range("M116").Value = 9.777799999999843
$ws.Range("H136").Value = 3983.8667
$ws.Range("I136").Value = 3519.923
$ws.Range("K136").Value = 10559.769
$ws.Range("M136").Value = -8009.769
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2305.8
$ws.Range("I3").Value = 2284.2222
$ws.Range("K3").Value = 2284.2222
$ws.Range("M3").Value = -2170.2222
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = ""
$ws.Range("H107").Value = 2627.4119
$ws.Range("I107").Value = 2296.6365
$ws.Range("J107").Value = 3233.8333
$ws.Range("K107").Value = 2296.6365
$ws.Range("L107").Value = 3233.8333
$ws.Range("M107").Value = -376.6365000000001
$ws.Range("N107").Value = -7073.8333
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2491.5
$ws.Range("I16").Value = 2212.375
$ws.Range("K16").Value = 2212.375
$ws.Range("M16").Value = -1925.375
$ws.Range("H31").Value = 45814.13
$ws.Range("I31").Value = 91709.28
$ws.Range("K31").Value = 91709.28
$ws.Range("M31").Value = -91414.28
$ws.Range("H34").Value = 45814.13
$ws.Range("I34").Value = 91709.28
$ws.Range("K34").Value = 91709.28
$ws.Range("M34").Value = -91507.28
$ws.Range("H58").Value = 3206.375
$ws.Range("I58").Value = 2550.2104
$ws.Range("J58").Value = 5699.8
$ws.Range("K58").Value = 2550.2104
$ws.Range("L58").Value = 5699.8
$ws.Range("M58").Value = -2347.2104
$ws.Range("N58").Value = -6105.8
$ws.Range("H107").Value = 718.3913
$ws.Range("I107").Value = 759.0952
$ws.Range("K107").Value = 759.0952
$ws.Range("M107").Value = 1160.9048
$ws.Range("H113").Value = 2491.5
$ws.Range("I113").Value = 2212.375
$ws.Range("K113").Value = 2212.375
$ws.Range("M113").Value = -42.375
$ws.Range("H136").Value = 3206.375
$ws.Range("I136").Value = 2550.2104
$ws.Range("J136").Value = 5699.8
$ws.Range("K136").Value = 7650.6312
$ws.Range("L136").Value = 17099.4
$ws.Range("M136").Value = -5100.6312
$ws.Range("N136").Value = -22199.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 524.2
$ws.Range("I23").Value = 375
$ws.Range("J23").Value = 561.5
$ws.Range("K23").Value = 1125
$ws.Range("L23").Value = 1684.5
$ws.Range("M23").Value = -890
$ws.Range("N23").Value = -2154.5
$ws.Range("H38").Value = 77.05882
$ws.Range("I38").Value = 76.53846
$ws.Range("K38").Value = 229.61538
$ws.Range("M38").Value = 117.38462
$ws.Range("H70").Value = 5226.6665
$ws.Range("I70").Value = 380
$ws.Range("J70").Value = 7650
$ws.Range("K70").Value = 1140
$ws.Range("L70").Value = 22950
$ws.Range("M70").Value = -825
$ws.Range("N70").Value = -23580
$ws.Range("H73").Value = 5226.6665
$ws.Range("I73").Value = 380
$ws.Range("J73").Value = 7650
$ws.Range("K73").Value = 1140
$ws.Range("L73").Value = 22950
$ws.Range("M73").Value = -48
$ws.Range("N73").Value = -25134
$ws.Range("H80").Value = 500
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 500
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H138").Value = 8362.471
$ws.Range("I138").Value = 5748.091
$ws.Range("J138").Value = 13155.5
$ws.Range("K138").Value = 17244.273
$ws.Range("L138").Value = 39466.5
$ws.Range("M138").Value = -12104.273
$ws.Range("N138").Value = -49746.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1444.0769
$ws.Range("I102").Value = 1397.75
$ws.Range("K102").Value = 1397.75
$ws.Range("M102").Value = 224.25
$ws.Range("H107").Value = 959.44446
$ws.Range("I107").Value = 898.1818
$ws.Range("J107").Value = 1055.7142
$ws.Range("K107").Value = 898.1818
$ws.Range("L107").Value = 1055.7142
$ws.Range("M107").Value = 1021.8182
$ws.Range("N107").Value = -4895.7142
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4770.8125
$ws.Range("I16").Value = 4918.4614
$ws.Range("J16").Value = 4131
$ws.Range("K16").Value = 4918.4614
$ws.Range("L16").Value = 4131
$ws.Range("M16").Value = -4748.4614
$ws.Range("N16").Value = -4471
$ws.Range("H32").Value = 1486.75
$ws.Range("I32").Value = 1486.75
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1486.75
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1169.75
$ws.Range("N32").Value = ""
$ws.Range("H134").Value = 49150
$ws.Range("J134").Value = 49150
$ws.Range("L134").Value = 49150
$ws.Range("N134").Value = -59290
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 19999
$ws.Range("J26").Value = 19999
$ws.Range("L26").Value = 19999
$ws.Range("N26").Value = -20585
$ws.Range("H70").Value = 39459.25
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 39459.25
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 39459.25
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -40089.25
$ws.Range("H73").Value = 39459.25
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 39459.25
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 39459.25
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -41643.25
